$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 566 (pushes the existing 566..669 down to 568..671)
$ws.Rows.Item(566).Insert()
$ws.Rows.Item(566).Insert()

# --- New row 566 ---
$ws.Range("A566").Value = 4
$ws.Range("B566").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C566").Value = "Los Lagos"
$ws.Range("D566").Value = 44995
$ws.Range("E566").Value = 10
$ws.Range("F566").Value = 100112006
$ws.Range("G566").Value = "Repollo"
$ws.Range("H566").Value = "Crespo record"
$ws.Range("I566").Value = "Primera"
$ws.Range("J566").Value = 600
$ws.Range("K566").Value = 1800
$ws.Range("L566").Value = 1800
$ws.Range("M566").Value = 1800
$ws.Range("N566").Value = "`$/unidad"
$ws.Range("O566").Value = "Región Metropolitana"
$ws.Range("P566").Value = 1800
$ws.Range("Q566").Value = 1
$ws.Range("R566").Value = "Hortaliza"

# --- New row 567 ---
$ws.Range("A567").Value = 4
$ws.Range("B567").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C567").Value = "Los Lagos"
$ws.Range("D567").Value = 44995
$ws.Range("E567").Value = 10
$ws.Range("F567").Value = 100112006
$ws.Range("G567").Value = "Repollo"
$ws.Range("H567").Value = "Crespo record"
$ws.Range("I567").Value = "Segunda"
$ws.Range("J567").Value = 600
$ws.Range("K567").Value = 1600
$ws.Range("L567").Value = 1600
$ws.Range("M567").Value = 1600
$ws.Range("N567").Value = "`$/unidad"
$ws.Range("O567").Value = "Región Metropolitana"
$ws.Range("P567").Value = 1600
$ws.Range("Q567").Value = 1
$ws.Range("R567").Value = "Hortaliza"

# --- New row 670 (appended after the shifted data, copy of old row 668) ---
$ws.Range("A670").Value = 4
$ws.Range("B670").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C670").Value = "Los Lagos"
$ws.Range("D670").Value = 44327
$ws.Range("E670").Value = 10
$ws.Range("F670").Value = 100112006
$ws.Range("G670").Value = "Repollo"
$ws.Range("H670").Value = "Crespo record"
$ws.Range("I670").Value = "Primera"
$ws.Range("J670").Value = 700
$ws.Range("K670").Value = 1000
$ws.Range("L670").Value = 1000
$ws.Range("M670").Value = 1000
$ws.Range("N670").Value = "`$/unidad"
$ws.Range("O670").Value = "Región Metropolitana"
$ws.Range("P670").Value = 1000
$ws.Range("Q670").Value = 1
$ws.Range("R670").Value = "Hortaliza"
$ws.Range("D670").NumberFormat = $ws.Range("D669").NumberFormat

# --- New row 671 (appended after the shifted data, copy of old row 669) ---
$ws.Range("A671").Value = 4
$ws.Range("B671").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C671").Value = "Los Lagos"
$ws.Range("D671").Value = 44327
$ws.Range("E671").Value = 10
$ws.Range("F671").Value = 100112006
$ws.Range("G671").Value = "Repollo"
$ws.Range("H671").Value = "Crespo record"
$ws.Range("I671").Value = "Primera"
$ws.Range("J671").Value = 700
$ws.Range("K671").Value = 1000
$ws.Range("L671").Value = 1000
$ws.Range("M671").Value = 1000
$ws.Range("N671").Value = "`$/unidad"
$ws.Range("O671").Value = "Región del Maule"
$ws.Range("P671").Value = 1000
$ws.Range("Q671").Value = 1
$ws.Range("R671").Value = "Hortaliza"
$ws.Range("D671").NumberFormat = $ws.Range("D669").NumberFormat
